# Applies the "add aluminum and brass anvil blocks" update to the
# "新增物品" (add items) worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows appended to the item table (columns A-D).
# Row 11 repeats the existing "laterite" entry (as in the source data),
# rows 12-18 introduce the new aluminum items and the hammer/anvil tools.
$newRows = @(
    @(10, "laterite",         "红土",   "Laterite"),
    @(11, "aluminum_ingot",   "铝锭",   "Aluminum Ingot"),
    @(12, "aluminum_sheet",   "铝板",   "Aluminum Sheet"),
    @(13, "aluminum_block",   "铝块",   "Block Of Aluminum"),
    @(14, "aluminum_nugget",  "铝粒",   "Aluminum Nugget"),
    @(15, "andesite_hammer",  "安山锤", "Andesite Hammer"),
    @(16, "brass_hammer",     "黄铜锤", "Brass Hammer"),
    @(17, "brass_anvil",      "黄铜砧", "Brass Anvil")
)

$r = 11
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Final selection on the sheet, as left by the author.
$ws.Range("B25").Select()

# Restore the workbook window size recorded by the author's Excel session.
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 1476
    $win.Height = 702.75
} catch {
}
